$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row below the current data (column A) and append the new ticker there.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "GRT-USD"
